$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)

# --- Update summary counts on "Roses form" sheet (sheet2) ---
$ws2.Range("B1").Value = "updated_search"
$ws2.Range("B5").Value = 1965
$ws2.Range("B7").Value = 81
$ws2.Range("B9").Value = 81

# --- Rebuild the exclusion/summary block (rows 14-36) ---
# Clear everything first so stale shared-string references to rows/labels
# that are being removed (e.g. "Excluded full texts (no PECO element)" and
# "Studies included in the map") get dropped entirely.
$ws2.Range("A14:C36").Clear()

$ws2.Range("A14").Value = "Excluded full texts (not about fauna)"
$ws2.Range("B14").Value = 2

$ws2.Range("A15").Value = "Excluded full texts (not about forest)"
$ws2.Range("B15").Value = 4

$ws2.Range("A16").Value = "Excluded full texts (no biodiversity outcome)"
$ws2.Range("B16").Value = 1

$ws2.Range("A17").Value = "Excluded full texts (no exposure)"
$ws2.Range("B17").Value = 27

$ws2.Range("A18").Value = "Excluded full texts (not field based)"
$ws2.Range("B18").Value = 1

$ws2.Range("A19").Value = "Excluded full texts (no comparison)"
$ws2.Range("B19").Value = 1

$ws2.Range("A20").Value = "Excluded full texts (no mean/median)"
$ws2.Range("B20").Value = 2

$ws2.Range("A21").Value = "Excluded full texts (duplicated data)"
$ws2.Range("B21").Value = 1

$ws2.Range("A22").Value = "Total excluded texts at full-text"
$ws2.Range("B22").Formula = "=SUM(B13:B21,B10:B11)"

$ws2.Range("A23").Value = "Articles in the review"
$ws2.Range("B23").Value = 38

$ws2.Range("A24").Value = "Studies in the review"
$ws2.Range("B24").Value = 38

# Update selection to reflect where the cursor ended up after editing
$ws2.Range("A25").Select()
